$wb = $excel.ActiveWorkbook

# --- Sheet "2025" (sheet1.xml) ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 44809.96542313504
$ws.Range("B2").Value = 19576.62650831837
$ws.Range("E2").Value = 135535.499643962
$ws.Range("G2").Value = 42315.16049510826
$ws.Range("H2").Value = 488787.6484426507
$ws.Range("I2").Value = 403424.0040594552
$ws.Range("N2").Value = 50337.28263146494
$ws.Range("O2").Value = 69397.29901820578

# --- Sheet "2030" (sheet2.xml) ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 29601.27271984311
$ws.Range("B2").Value = 101099.0316572796
$ws.Range("E2").Value = 66023.91524459935
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 203788.0418410577
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 18425.67154306719
$ws.Range("O2").Value = 34031.99607318347

# --- Sheet "2035" (sheet3.xml) ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 82423.16832641534
$ws.Range("B2").Value = 16911.72700928294
$ws.Range("E2").Value = 0
$ws.Range("I2").Value = 178689.2597542998
$ws.Range("M2").Value = 31702.65741071548
$ws.Range("N2").Value = 9464.444644938132
$ws.Range("O2").Value = 53633.84592823405
